$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = [double]"0"
$ws.Cells.Item(3, 2).Value = [double]"-4.440892098500626E-14"
$ws.Cells.Item(4, 2).Value = [double]"0"
$ws.Cells.Item(5, 2).Value = [double]"4.440892098500626E-14"
$ws.Cells.Item(6, 2).Value = [double]"0.8761527144223624"
$ws.Cells.Item(7, 2).Value = [double]"0.2498186593231866"
$ws.Cells.Item(8, 2).Value = [double]"-0.2073814777540428"
$ws.Cells.Item(9, 2).Value = [double]"0.3094428711141628"
$ws.Cells.Item(10, 2).Value = [double]"0.5223161956339206"
$ws.Cells.Item(11, 2).Value = [double]"0.3117781822009169"
$ws.Cells.Item(12, 2).Value = [double]"0.3783498544355668"
$ws.Cells.Item(13, 2).Value = [double]"0.6542703491021484"
$ws.Cells.Item(14, 2).Value = [double]"-0.6157632300240357"
$ws.Cells.Item(15, 2).Value = [double]"0.01566947406670405"
$ws.Cells.Item(16, 2).Value = [double]"-1.024794128387363"
$ws.Cells.Item(17, 2).Value = [double]"0.4283030634637974"
$ws.Cells.Item(18, 2).Value = [double]"0.4210175484930634"
$ws.Cells.Item(19, 2).Value = [double]"0.3309484344824476"
$ws.Cells.Item(20, 2).Value = [double]"0.03534272045342401"
$ws.Cells.Item(21, 2).Value = [double]"-1.144790596790379"
$ws.Cells.Item(22, 2).Value = [double]"-0.3737413844400406"
$ws.Cells.Item(23, 2).Value = [double]"0.3599799282585359"
$ws.Cells.Item(24, 2).Value = [double]"0.5766501347738604"
$ws.Cells.Item(25, 2).Value = [double]"-0.5004758034602208"
$ws.Cells.Item(26, 2).Value = [double]"-0.7750463390001627"
$ws.Cells.Item(27, 2).Value = [double]"0.2146387198177946"
$ws.Cells.Item(28, 2).Value = [double]"0.02351224639369764"
$ws.Cells.Item(29, 2).Value = [double]"-0.5673139978460418"
$ws.Cells.Item(30, 2).Value = [double]"0.03242225351332007"
$ws.Cells.Item(31, 2).Value = [double]"-0.3625064515619281"
$ws.Cells.Item(32, 2).Value = [double]"-0.5809412570459083"
$ws.Cells.Item(33, 2).Value = [double]"0.3008481533817164"
$ws.Cells.Item(34, 2).Value = [double]"-0.4092135626299287"
$ws.Cells.Item(35, 2).Value = [double]"-0.7785741838770672"
$ws.Cells.Item(36, 2).Value = [double]"-0.9959164522061803"
$ws.Cells.Item(37, 2).Value = [double]"-0.02772281113909703"
$ws.Cells.Item(38, 2).Value = [double]"0.565922940159902"
$ws.Cells.Item(39, 2).Value = [double]"0.266820966430148"
